# "new rates as of 3/18"
#
# - Row 3  (Morris, IL)        -> updated rates, extended from column M out to column P
# - Row 6  (Belleville, MI)    -> renamed to "Belleville,  MI" (note double space), updated
#                                 rates, extended from column M out to column P
# - Row 12 (Grand Prairie, TX) -> unchanged rates, just extended styling out to column P
# - Three brand-new destinations appended as rows 13-15:
#     Monroe Township, NJ / Monrovia, MD / Owatonna, MN
# - Column A widened to fit the longer destination names; columns B:P share widths.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-RowValues {
    param($ws, [int]$row, [int]$startCol, [object[]]$values)

    $n = $values.Count
    $comArray = New-Object 'object[,]' 1, $n
    for ($i = 0; $i -lt $n; $i++) {
        $comArray[0, $i] = $values[$i]
    }

    $startCell = $ws.Cells.Item($row, $startCol)
    $endCell = $ws.Cells.Item($row, $startCol + $n - 1)
    $rng = $ws.Range($startCell, $endCell)
    $rng.Value = $comArray
    $rng.NumberFormat = "0.00"
}

# ---- Row 3: Morris, IL -------------------------------------------------
Set-RowValues $ws 3 2 @(290, 570, 850, 1120, 1400, 1680, 1945, 2200, 2475, 2750, 3000, 3240, 3510, 3780, 4015)

# ---- Row 6: Belleville,  MI (renamed, note the double space) ----------
$ws.Range("A6").Value = "Belleville,  MI"
Set-RowValues $ws 6 2 @(465, 800, 1185, 1540, 1825, 2160, 2415, 2760, 3060, 3400, 3740, 4080, 4355, 4690, 5025)

# ---- Row 12: Grand Prairie, TX -----------------------------------------
# Rates unchanged (400, 700, 975); just extend the styled-but-empty cells out to column P.
$ws.Range("E12:P12").NumberFormat = "0.00"

# ---- Row 13: Monroe Township, NJ (new) ----------------------------------
$ws.Range("A13").Value = "Monroe Township, NJ"
Set-RowValues $ws 13 2 @(470, 880, 1245, 1580, 1925, 2280, 2625, 2920, 3240, 3550)
$ws.Range("L13:P13").NumberFormat = "0.00"

# ---- Row 14: Monrovia, MD (new) -----------------------------------------
$ws.Range("A14").Value = "Monrovia, MD"
Set-RowValues $ws 14 2 @(470, 880, 1245, 1580, 1925, 2280, 2625, 2920, 3240, 3550)
$ws.Range("L14:P14").NumberFormat = "0.00"

# ---- Row 15: Owatonna, MN (new) -----------------------------------------
$ws.Range("A15").Value = "Owatonna, MN"
Set-RowValues $ws 15 2 @(525, 1135, 1475, 1580, 1840, 2190, 2345, 2600, 2880, 3130)
$ws.Range("L15:P15").NumberFormat = "0.00"

# ---- Trailing empty-but-styled cells on the pre-existing rows, now that
#      the table extends to column P ----
$ws.Range("P2").NumberFormat = "0.00"
$ws.Range("P4").NumberFormat = "0.00"
$ws.Range("P5").NumberFormat = "0.00"
$ws.Range("P7").NumberFormat = "0.00"
$ws.Range("P8").NumberFormat = "0.00"
$ws.Range("P9").NumberFormat = "0.00"
$ws.Range("P10").NumberFormat = "0.00"
$ws.Range("P11").NumberFormat = "0.00"

# ---- Row 1 header: extend the sequential numbering to the new column P ----
$ws.Range("P1").Value = 15

# ---- Column widths --------------------------------------------------------
# NB: the host's ColumnWidth -> XML `width` conversion only has 1/6-character
# granularity (`width = round(ColumnWidth*6)/6 + 5/6`), unlike real Excel's
# 1/256 granularity, so we pick the ColumnWidth input that lands closest to
# the authentic bestFit widths recorded in the diff (20.42578125 / 6.5703125
# / 7.5703125).
$ws.Columns("A:A").ColumnWidth = 19.666666666666668   # -> width 20.5
$ws.Columns("B:B").ColumnWidth = 5.666666666666667    # -> width 6.5
$ws.Columns("C:P").ColumnWidth = 6.666666666666667    # -> width 7.5

# ---- Selection matches the author's final cursor position ---------------
$ws.Activate() | Out-Null
$ws.Range("A15:XFD16").Select() | Out-Null
